# Auto-generated Excel COM-interop script
# Applies numeric updates to the LeveProfit/averagePrice columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4507.294
$ws.Range("I74").Value = 5871.875
$ws.Range("J74").Value = 3294.3333
$ws.Range("K74").Value = 5871.875
$ws.Range("L74").Value = 3294.3333
$ws.Range("M74").Value = -4935.875
$ws.Range("N74").Value = -5166.3333

$ws.Range("H77").Value = 4507.294
$ws.Range("I77").Value = 5871.875
$ws.Range("J77").Value = 3294.3333
$ws.Range("K77").Value = 29359.375
$ws.Range("L77").Value = 16471.6665
$ws.Range("M77").Value = -24679.375
$ws.Range("N77").Value = -25831.6665

$ws.Range("H92").Value = 1745
$ws.Range("I92").Value = 2226.6667
$ws.Range("J92").Value = 878
$ws.Range("K92").Value = 2226.6667
$ws.Range("L92").Value = 878
$ws.Range("M92").Value = -978.6667000000002
$ws.Range("N92").Value = -3374

$ws.Range("H103").Value = 1950.0834
$ws.Range("I103").Value = 2330
$ws.Range("J103").Value = 1678.7142
$ws.Range("K103").Value = 6990
$ws.Range("L103").Value = 5036.142599999999
$ws.Range("M103").Value = -6404
$ws.Range("N103").Value = -6208.142599999999

$ws.Range("H113").Value = 335001.66
$ws.Range("I113").Value = 335001.66
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 335001.66
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -331747.66
$ws.Range("N113").ClearContents()

$ws.Range("H129").Value = 2835.0588
$ws.Range("J129").Value = 942.81396
$ws.Range("L129").Value = 2828.44188
$ws.Range("N129").Value = -12828.44188

$ws.Range("H138").Value = 4047.4312
$ws.Range("I138").Value = 3389.1
$ws.Range("J138").Value = 4184.5835
$ws.Range("K138").Value = 10167.3
$ws.Range("L138").Value = 12553.7505
$ws.Range("M138").Value = -5027.299999999999
$ws.Range("N138").Value = -22833.7505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 77995.46000000001
$ws.Range("I2").Value = 1161.75
$ws.Range("J2").Value = 1000000
$ws.Range("K2").Value = 1161.75
$ws.Range("L2").Value = 1000000
$ws.Range("M2").Value = -1048.75
$ws.Range("N2").Value = -1000226

$ws.Range("H32").Value = 6747.73
$ws.Range("I32").Value = 6477.2754
$ws.Range("K32").Value = 6477.2754
$ws.Range("M32").Value = -6190.2754

$ws.Range("H63").Value = 3347.75
$ws.Range("J63").Value = 3996.6667
$ws.Range("L63").Value = 3996.6667
$ws.Range("N63").Value = -5368.6667

$ws.Range("H66").Value = 3347.75
$ws.Range("J66").Value = 3996.6667
$ws.Range("L66").Value = 19983.3335
$ws.Range("N66").Value = -26847.3335

$ws.Range("H116").Value = 77995.46000000001
$ws.Range("I116").Value = 1161.75
$ws.Range("J116").Value = 1000000
$ws.Range("K116").Value = 1161.75
$ws.Range("L116").Value = 1000000
$ws.Range("M116").Value = 1132.25
$ws.Range("N116").Value = -1004588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 77995.46000000001
$ws.Range("I3").Value = 1161.75
$ws.Range("J3").Value = 1000000
$ws.Range("K3").Value = 1161.75
$ws.Range("L3").Value = 1000000
$ws.Range("M3").Value = -1047.75
$ws.Range("N3").Value = -1000228

$ws.Range("H82").Value = 15151
$ws.Range("J82").Value = 33131.5
$ws.Range("L82").Value = 33131.5
$ws.Range("N82").Value = -33897.5

$ws.Range("H85").Value = 15151
$ws.Range("J85").Value = 33131.5
$ws.Range("L85").Value = 33131.5
$ws.Range("N85").Value = -35783.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 123214.75
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 123214.75
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 123214.75
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -123804.75

$ws.Range("H34").Value = 123214.75
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 123214.75
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 123214.75
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -123618.75

$ws.Range("H68").Value = 17020.055
$ws.Range("J68").Value = 17020.055
$ws.Range("L68").Value = 17020.055
$ws.Range("N68").Value = -18518.055

$ws.Range("H71").Value = 17020.055
$ws.Range("J71").Value = 17020.055
$ws.Range("L71").Value = 51060.165
$ws.Range("N71").Value = -58548.165

$ws.Range("H74").Value = 38712
$ws.Range("J74").Value = 38712
$ws.Range("L74").Value = 38712
$ws.Range("N74").Value = -40460

$ws.Range("H77").Value = 38712
$ws.Range("J77").Value = 38712
$ws.Range("L77").Value = 116136
$ws.Range("N77").Value = -124872

$ws.Range("H92").Value = 100000
$ws.Range("J92").Value = 100000
$ws.Range("L92").Value = 100000
$ws.Range("N92").Value = -104992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 864.4286
$ws.Range("I4").Value = 275.5
$ws.Range("J4").Value = 1100
$ws.Range("K4").Value = 826.5
$ws.Range("L4").Value = 3300
$ws.Range("M4").Value = -714.5
$ws.Range("N4").Value = -3524

$ws.Range("H7").Value = 559.5
$ws.Range("I7").Value = 359
$ws.Range("J7").Value = 760
$ws.Range("K7").Value = 1077
$ws.Range("L7").Value = 2280
$ws.Range("M7").Value = -965
$ws.Range("N7").Value = -2504

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2750
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 4500
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 4500
$ws.Range("M5").Value = -888
$ws.Range("N5").Value = -4724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 365636.38
$ws.Range("J2").Value = 9000
$ws.Range("L2").Value = 9000
$ws.Range("N2").Value = -9224

$ws.Range("H136").Value = 2607.4517
$ws.Range("I136").Value = 1745.6666
$ws.Range("K136").Value = 5236.9998
$ws.Range("M136").Value = -2686.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 17641
$ws.Range("I2").Value = 7996.3335
$ws.Range("J2").Value = 24874.5
$ws.Range("K2").Value = 7996.3335
$ws.Range("L2").Value = 24874.5
$ws.Range("M2").Value = -7884.3335
$ws.Range("N2").Value = -25098.5

$ws.Range("H107").Value = 143445.14
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 167186
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 501558
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -505398

$ws.Range("H133").Value = 36666.668
$ws.Range("J133").Value = 36666.668
$ws.Range("L133").Value = 36666.668
$ws.Range("N133").Value = -46786.668
